# PR #32: append a new data row (issue #18 / "テスト") to the tracking
# sheet, growing the used range from A1:D1 to A1:D2.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "#18"
$ws.Range("B2").Value = "テスト"
